$wb = $excel.ActiveWorkbook

# Rename the metadata sheet
$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Make the renamed sheet the active sheet/tab
$metaSheet.Activate()
